$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Fill in the detail hours for each person/role. Cells that stay at zero
# are written as the literal placeholder text "-" (as the author did),
# everything else gets its real hour count.
# ---------------------------------------------------------------------

# Row 2 - Andrea Favero
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"

# Row 3 - Eleonora Thiella
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3

# Row 4 - Federico Caldart
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 3

# Row 5 - Giovanni Cavallin
$ws.Range("B5").Value = "-"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 5

# Row 6 - Giovanni Dalla Riva
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 2

# Row 7 - Lorenzo Menegon
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = "-"
$ws.Range("G7").Value = "-"

# Row 8 - Stefano Panozzo
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = 2

# The old hard-coded helper totals in row 11 are no longer needed now
# that row 9 computes the real totals from the data above - clear them
# but keep their styling.
$ws.Range("B11:G11").ClearContents()

# Update the active selection/view to match where the author left off.
$ws.Range("A1:H9").Select()
$ws.Range("H9").Activate()

# ---------------------------------------------------------------------
# Move/resize the chart so it sits to the right of the table, spanning
# columns I..U (0-based 8..20) and rows 1..10 (0-based 0..9).
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 690.902283464567
$co.Top = 0.0
$co.Width = 736.8749606299212
$co.Height = 287.28748031496065
